$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix class schedule entries for turm ELT-2A.
# The previous data erroneously showed only a single class; restore the
# correct set of entries so multiple classes are reflected.

$ws.Range("F2").Value = "-"
$ws.Range("C3").Value = "Circuitos Elétricos 2"
$ws.Range("F9").Value = "Circuitos Elétricos 2"
$ws.Range("B10").Value = "-"
$ws.Range("F10").Value = "Circuitos Elétricos 2"
$ws.Range("F11").Value = "-"
$ws.Range("B13").Value = "-"
